# chore: update Sheets via scheduled runner
# Refreshes cached marketboard price / profit figures (columns H-N:
# currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
# LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) for specific
# leve rows across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
# Row 12
$ws.Range("H12").Value = 446.83334
$ws.Range("I12").Value = 446.83334
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 446.83334
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -276.83334
# Row 17
$ws.Range("H17").Value = 5835.5
$ws.Range("J17").Value = 5835.5
$ws.Range("L17").Value = 17506.5
$ws.Range("N17").Value = -17842.5
# Row 33
$ws.Range("H33").Value = 1594.7273
$ws.Range("I33").Value = 388.22223
$ws.Range("K33").Value = 388.22223
$ws.Range("M33").Value = -159.22223
# Row 62
$ws.Range("H62").Value = 5569.4346
$ws.Range("I62").Value = 4599.8335
$ws.Range("K62").Value = 4599.8335
$ws.Range("M62").Value = -3975.8335
# Row 65
$ws.Range("H65").Value = 5569.4346
$ws.Range("I65").Value = 4599.8335
$ws.Range("K65").Value = 22999.1675
$ws.Range("M65").Value = -19879.1675

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Range("H2").Value = 458.58334
$ws.Range("I2").Value = 439.09525
$ws.Range("K2").Value = 439.09525
$ws.Range("M2").Value = -326.09525
# Row 32
$ws.Range("H32").Value = 12019.84
$ws.Range("I32").Value = 8749.333000000001
$ws.Range("K32").Value = 8749.333000000001
$ws.Range("M32").Value = -8462.333000000001
# Row 45
$ws.Range("H45").Value = 3051.7144
$ws.Range("I45").Value = 1589.7142
$ws.Range("J45").Value = 4513.7144
$ws.Range("K45").Value = 1589.7142
$ws.Range("L45").Value = 4513.7144
$ws.Range("M45").Value = -1212.7142
$ws.Range("N45").Value = -5267.7144
# Row 63
$ws.Range("H63").Value = 2363.2856
$ws.Range("I63").Value = 2164.6667
$ws.Range("K63").Value = 2164.6667
$ws.Range("M63").Value = -1478.6667
# Row 66
$ws.Range("H66").Value = 2363.2856
$ws.Range("I66").Value = 2164.6667
$ws.Range("K66").Value = 10823.3335
$ws.Range("M66").Value = -7391.333500000001
# Row 116
$ws.Range("H116").Value = 458.58334
$ws.Range("I116").Value = 439.09525
$ws.Range("K116").Value = 439.09525
$ws.Range("M116").Value = 1854.90475
# Row 124
$ws.Range("H124").Value = 37499.5
$ws.Range("J124").Value = 37499.5
$ws.Range("L124").Value = 37499.5
$ws.Range("N124").Value = -47319.5
# Row 125
$ws.Range("H125").Value = 75905
$ws.Range("J125").Value = 75905
$ws.Range("L125").Value = 75905
$ws.Range("N125").Value = -85745

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Range("H3").Value = 458.58334
$ws.Range("I3").Value = 439.09525
$ws.Range("K3").Value = 439.09525
$ws.Range("M3").Value = -325.09525
# Row 22
$ws.Range("H22").Value = 1020.9259
$ws.Range("I22").Value = 921.2222
$ws.Range("K22").Value = 921.2222
$ws.Range("M22").Value = -748.2222
# Row 119
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
# Row 134
$ws.Range("H134").Value = 3060.4583
$ws.Range("I134").Value = 2410.913
$ws.Range("K134").Value = 7232.739
$ws.Range("M134").Value = -4697.739

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
# Row 7
$ws.Range("H7").Value = 1613.3235
$ws.Range("I7").Value = 695.5
$ws.Range("J7").Value = 3296
$ws.Range("K7").Value = 695.5
$ws.Range("L7").Value = 3296
$ws.Range("M7").Value = -582.5
$ws.Range("N7").Value = -3522
# Row 31
$ws.Range("H31").Value = 7592.8076
$ws.Range("I31").Value = 4634.364
$ws.Range("J31").Value = 12731.158
$ws.Range("K31").Value = 4634.364
$ws.Range("L31").Value = 12731.158
$ws.Range("M31").Value = -4339.364
$ws.Range("N31").Value = -13321.158
# Row 34
$ws.Range("H34").Value = 7592.8076
$ws.Range("I34").Value = 4634.364
$ws.Range("J34").Value = 12731.158
$ws.Range("K34").Value = 4634.364
$ws.Range("L34").Value = 12731.158
$ws.Range("M34").Value = -4432.364
$ws.Range("N34").Value = -13135.158

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
# Row 112
$ws.Range("H112").Value = 14748.625
$ws.Range("I112").Value = 9329.666999999999
$ws.Range("K112").Value = 27989.001
$ws.Range("M112").Value = -26881.001
# Row 129
$ws.Range("H129").Value = 3222.95
$ws.Range("J129").Value = 5558.8887
$ws.Range("L129").Value = 16676.6661
$ws.Range("N129").Value = -26676.6661
# Row 136
$ws.Range("H136").Value = 1906
$ws.Range("I136").Value = 1906
$ws.Range("K136").Value = 5718
$ws.Range("M136").Value = -618
# Row 139
$ws.Range("H139").Value = 2374.1333
$ws.Range("J139").Value = 2500
$ws.Range("L139").Value = 7500
$ws.Range("N139").Value = -17780

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
# Row 2
$ws.Range("H2").Value = 1137.3684
$ws.Range("I2").Value = 1383.9231
$ws.Range("J2").Value = 603.1667
$ws.Range("K2").Value = 1383.9231
$ws.Range("L2").Value = 603.1667
$ws.Range("M2").Value = -1270.9231
$ws.Range("N2").Value = -829.1667
# Row 102
$ws.Range("H102").Value = 2007.3235
$ws.Range("I102").Value = 1397.5555
$ws.Range("J102").Value = 2693.3125
$ws.Range("K102").Value = 1397.5555
$ws.Range("L102").Value = 2693.3125
$ws.Range("M102").Value = 224.4445000000001
$ws.Range("N102").Value = -5937.3125
# Row 113
$ws.Range("H113").Value = 2009.9445
$ws.Range("I113").Value = 1774.4375
$ws.Range("K113").Value = 1774.4375
$ws.Range("M113").Value = 395.5625
# Row 126
$ws.Range("H126").Value = 10875.333
$ws.Range("I126").Value = 14856
$ws.Range("J126").Value = 2914
$ws.Range("K126").Value = 44568
$ws.Range("L126").Value = 8742
$ws.Range("M126").Value = -42098
$ws.Range("N126").Value = -13682
# Row 132
$ws.Range("H132").Value = 4418.9585
$ws.Range("I132").Value = 3325.4814
$ws.Range("J132").Value = 5824.857
$ws.Range("K132").Value = 9976.4442
$ws.Range("L132").Value = 17474.571
$ws.Range("M132").Value = -7446.4442
$ws.Range("N132").Value = -22534.571

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
# Row 22
$ws.Range("H22").Value = 2944.3333
$ws.Range("I22").Value = 2383.875
$ws.Range("J22").Value = 3584.8572
$ws.Range("K22").Value = 2383.875
$ws.Range("L22").Value = 3584.8572
$ws.Range("M22").Value = -2088.875
$ws.Range("N22").Value = -4174.8572
# Row 27
$ws.Range("H27").Value = 2944.3333
$ws.Range("I27").Value = 2383.875
$ws.Range("J27").Value = 3584.8572
$ws.Range("K27").Value = 2383.875
$ws.Range("L27").Value = 3584.8572
$ws.Range("M27").Value = -2276.875
$ws.Range("N27").Value = -3798.8572
# Row 40
$ws.Range("H40").Value = 6802.5557
$ws.Range("I40").Value = 4300.6
$ws.Range("K40").Value = 4300.6
$ws.Range("M40").Value = -4164.6
# Row 46
$ws.Range("H46").Value = 3507.6667
$ws.Range("I46").Value = 2060.25
$ws.Range("K46").Value = 2060.25
$ws.Range("M46").Value = -1872.25
# Row 99
$ws.Range("H99").Value = 45313.6
$ws.Range("J99").Value = 69284.5
$ws.Range("L99").Value = 69284.5
$ws.Range("N99").Value = -75274.5
# Row 122
$ws.Range("H122").Value = 9793.066000000001
$ws.Range("I122").Value = 9039.1
$ws.Range("K122").Value = 27117.3
$ws.Range("M122").Value = -24667.3

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
# Row 96
$ws.Range("H96").Value = 5598.8887
$ws.Range("I96").Value = 2198.9092
$ws.Range("J96").Value = 10941.714
$ws.Range("K96").Value = 2198.9092
$ws.Range("L96").Value = 10941.714
$ws.Range("M96").Value = -825.9092000000001
$ws.Range("N96").Value = -13687.714
# Row 100
$ws.Range("H100").Value = 717.8261
$ws.Range("I100").Value = 339.57144
$ws.Range("J100").Value = 1306.2222
$ws.Range("K100").Value = 679.14288
$ws.Range("L100").Value = 2612.4444
$ws.Range("M100").Value = -138.14288
$ws.Range("N100").Value = -3694.4444
# Row 107
$ws.Range("H107").Value = 1272.234
$ws.Range("J107").Value = 1009.125
$ws.Range("L107").Value = 3027.375
$ws.Range("N107").Value = -6867.375
# Row 132
$ws.Range("H132").Value = 4658.304
$ws.Range("I132").Value = 3530.2856
$ws.Range("J132").Value = 16502.5
$ws.Range("K132").Value = 10590.8568
$ws.Range("L132").Value = 49507.5
$ws.Range("M132").Value = -8060.856800000001
$ws.Range("N132").Value = -54567.5
